$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 9561
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = $null
$ws.Range("H23").Value = 9561
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = $null
$ws.Range("H69").Value = 7749.4443
$ws.Range("H70").Value = 2466.6667
$ws.Range("I70").Value = 2400
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 7200
$ws.Range("L70").Value = 7500
$ws.Range("M70").Value = -6930
$ws.Range("N70").Value = -8040
$ws.Range("H72").Value = 7749.4443
$ws.Range("H73").Value = 2466.6667
$ws.Range("I73").Value = 2400
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 7200
$ws.Range("L73").Value = 7500
$ws.Range("M73").Value = -6264
$ws.Range("N73").Value = -9372
$ws.Range("H86").Value = 1479.6
$ws.Range("I86").Value = 1599.5
$ws.Range("K86").Value = 1599.5
$ws.Range("M86").Value = -476.5
$ws.Range("H89").Value = 1479.6
$ws.Range("I89").Value = 1599.5
$ws.Range("K89").Value = 7997.5
$ws.Range("M89").Value = -2381.5
$ws.Range("H116").Value = 5790.5386
$ws.Range("I116").Value = 5762.9
$ws.Range("J116").Value = 5882.6665
$ws.Range("K116").Value = 5762.9
$ws.Range("L116").Value = 5882.6665
$ws.Range("M116").Value = -2320.9
$ws.Range("N116").Value = -12766.6665
$ws.Range("H135").Value = 621.7646999999999
$ws.Range("I135").Value = 438.13333
$ws.Range("K135").Value = 3943.19997
$ws.Range("M135").Value = -1408.19997
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5247.75
$ws.Range("I88").Value = 2991
$ws.Range("K88").Value = 2991
$ws.Range("M88").Value = -2585
$ws.Range("H91").Value = 5247.75
$ws.Range("I91").Value = 2991
$ws.Range("K91").Value = 2991
$ws.Range("M91").Value = -1587
$ws.Range("H102").Value = 5225.6
$ws.Range("I102").Value = 4532
$ws.Range("K102").Value = 4532
$ws.Range("M102").Value = -2910
$ws.Range("H110").Value = 1547.0588
$ws.Range("I110").Value = 1475.5
$ws.Range("K110").Value = 1475.5
$ws.Range("M110").Value = 569.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2306.25
$ws.Range("I20").Value = 2370
$ws.Range("K20").Value = 2370
$ws.Range("M20").Value = -2123
$ws.Range("H94").Value = 6615
$ws.Range("I94").Value = 961
$ws.Range("J94").Value = 10384.333
$ws.Range("K94").Value = 961
$ws.Range("L94").Value = 10384.333
$ws.Range("M94").Value = -510
$ws.Range("N94").Value = -11286.333
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 854.6667
$ws.Range("I16").Value = 854.6667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 854.6667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -567.6667
$ws.Range("N16").Value = $null
$ws.Range("H31").Value = 2695
$ws.Range("I31").Value = 1993.3572
$ws.Range("K31").Value = 1993.3572
$ws.Range("M31").Value = -1698.3572
$ws.Range("H34").Value = 2695
$ws.Range("I34").Value = 1993.3572
$ws.Range("K34").Value = 1993.3572
$ws.Range("M34").Value = -1791.3572
$ws.Range("H107").Value = 385.5
$ws.Range("I107").Value = 371.55554
$ws.Range("K107").Value = 371.55554
$ws.Range("M107").Value = 1548.44446
$ws.Range("H113").Value = 854.6667
$ws.Range("I113").Value = 854.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 854.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1315.3333
$ws.Range("N113").Value = $null
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 554.5
$ws.Range("J17").Value = 859.6667
$ws.Range("L17").Value = 2579.0001
$ws.Range("N17").Value = -2917.0001
$ws.Range("H34").Value = 2897.5454
$ws.Range("J34").Value = 4595.769
$ws.Range("L34").Value = 13787.307
$ws.Range("N34").Value = -13955.307
$ws.Range("H39").Value = 3034.875
$ws.Range("J39").Value = 3034.875
$ws.Range("L39").Value = 9104.625
$ws.Range("N39").Value = -9692.625
$ws.Range("H55").Value = 1924.875
$ws.Range("J55").Value = 2925
$ws.Range("L55").Value = 8775
$ws.Range("N55").Value = -9129
$ws.Range("H56").Value = 3528.963
$ws.Range("I56").Value = 3528.963
$ws.Range("K56").Value = 3528.963
$ws.Range("M56").Value = -2998.963
$ws.Range("H121").Value = 1420809.8
$ws.Range("I121").Value = 214458.86
$ws.Range("K121").Value = 643376.58
$ws.Range("M121").Value = -642066.58
$ws.Range("H126").Value = 2030
$ws.Range("I126").Value = 2030
$ws.Range("K126").Value = 6090
$ws.Range("M126").Value = -1150
$ws.Range("H132").Value = 3059.6
$ws.Range("J132").Value = 1500
$ws.Range("L132").Value = 13500
$ws.Range("N132").Value = -18560
$ws.Range("H133").Value = 5069.4287
$ws.Range("I133").Value = 3997
$ws.Range("J133").Value = 5873.75
$ws.Range("K133").Value = 11991
$ws.Range("L133").Value = 17621.25
$ws.Range("M133").Value = -6931
$ws.Range("N133").Value = -27741.25
$ws.Range("H134").Value = 9497
$ws.Range("I134").Value = 9497
$ws.Range("K134").Value = 28491
$ws.Range("M134").Value = -23421
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = $null
$ws.Range("H46").Value = 3279.875
$ws.Range("I46").Value = 3279.875
$ws.Range("K46").Value = 3279.875
$ws.Range("M46").Value = -3123.875
$ws.Range("H57").Value = 17333.334
$ws.Range("H70").Value = 7313.3335
$ws.Range("I70").Value = 3660
$ws.Range("K70").Value = 3660
$ws.Range("M70").Value = -3390
$ws.Range("H73").Value = 7313.3335
$ws.Range("I73").Value = 3660
$ws.Range("K73").Value = 3660
$ws.Range("M73").Value = -2724
$ws.Range("H80").Value = 6079.3477
$ws.Range("I80").Value = 3885.3333
$ws.Range("J80").Value = 8472.817999999999
$ws.Range("K80").Value = 3885.3333
$ws.Range("L80").Value = 8472.817999999999
$ws.Range("M80").Value = -2887.3333
$ws.Range("N80").Value = -10468.818
$ws.Range("H83").Value = 6079.3477
$ws.Range("I83").Value = 3885.3333
$ws.Range("J83").Value = 8472.817999999999
$ws.Range("K83").Value = 19426.6665
$ws.Range("L83").Value = 42364.09
$ws.Range("M83").Value = -14434.6665
$ws.Range("N83").Value = -52348.09
$ws.Range("H97").Value = 878.5714
$ws.Range("I97").Value = 871
$ws.Range("J97").Value = 897.5
$ws.Range("K97").Value = 871
$ws.Range("L97").Value = 897.5
$ws.Range("M97").Value = -375
$ws.Range("N97").Value = -1889.5
$ws.Range("H102").Value = 3518
$ws.Range("I102").Value = 3675
$ws.Range("J102").Value = 2890
$ws.Range("K102").Value = 3675
$ws.Range("L102").Value = 2890
$ws.Range("M102").Value = -2053
$ws.Range("N102").Value = -6134
$ws.Range("H107").Value = 1422.0834
$ws.Range("I107").Value = 786
$ws.Range("K107").Value = 786
$ws.Range("M107").Value = 1134
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3781.4285
$ws.Range("I22").Value = 2944.4285
$ws.Range("J22").Value = 4618.4287
$ws.Range("K22").Value = 2944.4285
$ws.Range("L22").Value = 4618.4287
$ws.Range("M22").Value = -2649.4285
$ws.Range("N22").Value = -5208.4287
$ws.Range("H27").Value = 3781.4285
$ws.Range("I27").Value = 2944.4285
$ws.Range("J27").Value = 4618.4287
$ws.Range("K27").Value = 2944.4285
$ws.Range("L27").Value = 4618.4287
$ws.Range("M27").Value = -2837.4285
$ws.Range("N27").Value = -4832.4287
$ws.Range("H40").Value = 4109.1763
$ws.Range("I40").Value = 4068
$ws.Range("J40").Value = 4138
$ws.Range("K40").Value = 4068
$ws.Range("L40").Value = 4138
$ws.Range("M40").Value = -3932
$ws.Range("N40").Value = -4410
$ws.Range("H82").Value = 2212
$ws.Range("I82").Value = 2301
$ws.Range("K82").Value = 2301
$ws.Range("M82").Value = -1940
$ws.Range("H85").Value = 2212
$ws.Range("I85").Value = 2301
$ws.Range("K85").Value = 2301
$ws.Range("M85").Value = -1053
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 57686
$ws.Range("I63").Value = 39393
$ws.Range("J63").Value = 75979
$ws.Range("K63").Value = 39393
$ws.Range("L63").Value = 75979
$ws.Range("M63").Value = -38769
$ws.Range("N63").Value = -77227
$ws.Range("H66").Value = 57686
$ws.Range("I66").Value = 39393
$ws.Range("J66").Value = 75979
$ws.Range("K66").Value = 118179
$ws.Range("L66").Value = 227937
$ws.Range("M66").Value = -115059
$ws.Range("N66").Value = -234177
